$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Portfolio Name" header in column L, following the existing header row
$ws.Range("L1").Value = "Portfolio Name"

# Update the active selection to reflect the next empty cell after the new column
$ws.Range("M2").Select()
